# Applies the updated cryptocurrency price/volume figures, and the
# corrected row order for BitcoinCash / ImmutableX (rows 23-24),
# as published by the scheduled GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.538.15"
$ws.Range("E2").Value = "  +4.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.472.45"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.95"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.78"
$ws.Range("E6").Value = "  +4.20%  "
$ws.Range("E7").Value = "  +1.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.37"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.865.91"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.480.69"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.847"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.453.09"
$ws.Range("E18").Value = "  +4.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("E19").Value = "  +2.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.48"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.59"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.69"
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  +3.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("E25").Value = "  +2.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.17"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  -4.01%  "
$ws.Range("E29").Value = "  +2.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.72"
$ws.Range("E30").Value = "  +3.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.70"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.75"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.95"
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.53"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.76"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.984.98"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.82"
$ws.Range("E48").Value = "  +8.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.07"
$ws.Range("E49").Value = "  -3.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.22"
$ws.Range("E50").Value = "  +12.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.15"
$ws.Range("E51").Value = "  +5.38%  "
